# Applies the cell-value updates described by the commit diff
# (Leve profit-tracking sheets refreshed from the market-board scrape).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 446
$ws.Range("I2").Value = 325.2
$ws.Range("J2").Value = 618.5714
$ws.Range("K2").Value = 325.2
$ws.Range("L2").Value = 618.5714
$ws.Range("M2").Value = -212.2
$ws.Range("N2").Value = -844.5714
$ws.Range("H9").Value = 116.8125
$ws.Range("I9").Value = 119
$ws.Range("J9").Value = 112
$ws.Range("K9").Value = 119
$ws.Range("L9").Value = 112
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = -450
$ws.Range("H38").Value = 204.28572
$ws.Range("J38").Value = 333
$ws.Range("L38").Value = 999
$ws.Range("N38").Value = -1743
$ws.Range("H51").Value = 2438.4614
$ws.Range("J51").Value = 2433.3333
$ws.Range("L51").Value = 2433.3333
$ws.Range("N51").Value = -3401.3333
$ws.Range("H58").Value = 1672.6842
$ws.Range("I58").Value = 1202.3334
$ws.Range("J58").Value = 1889.7693
$ws.Range("K58").Value = 3607.0002
$ws.Range("L58").Value = 5669.3079
$ws.Range("M58").Value = -3457.0002
$ws.Range("N58").Value = -5969.3079
$ws.Range("H87").Value = 26135.4
$ws.Range("J87").Value = 26135.4
$ws.Range("L87").Value = 26135.4
$ws.Range("N87").Value = -28631.4
$ws.Range("H90").Value = 26135.4
$ws.Range("J90").Value = 26135.4
$ws.Range("L90").Value = 78406.20000000001
$ws.Range("N90").Value = -90886.20000000001
$ws.Range("H138").Value = 3331.1428
$ws.Range("I138").Value = 2286.6667
$ws.Range("J138").Value = 3748.9333
$ws.Range("K138").Value = 6860.000100000001
$ws.Range("L138").Value = 11246.7999
$ws.Range("M138").Value = -1720.000100000001
$ws.Range("N138").Value = -21526.7999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 58259.43
$ws.Range("J130").Value = 58259.43
$ws.Range("L130").Value = 58259.43
$ws.Range("N130").Value = -68299.42999999999
$ws.Range("H132").Value = 2393.6428
$ws.Range("I132").Value = 1927.6765
$ws.Range("J132").Value = 4374
$ws.Range("K132").Value = 5783.029500000001
$ws.Range("L132").Value = 13122
$ws.Range("M132").Value = -3253.029500000001
$ws.Range("N132").Value = -18182

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1254318.2
$ws.Range("I105").Value = 1842633
$ws.Range("J105").Value = 4149.625
$ws.Range("K105").Value = 1842633
$ws.Range("L105").Value = 4149.625
$ws.Range("M105").Value = -1840886
$ws.Range("N105").Value = -7643.625
$ws.Range("H131").Value = 27950
$ws.Range("J131").Value = 27950
$ws.Range("L131").Value = 27950
$ws.Range("N131").Value = -38030
$ws.Range("H132").Value = 50599.5
$ws.Range("J132").Value = 50599.5
$ws.Range("L132").Value = 50599.5
$ws.Range("N132").Value = -60719.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 241.42857
$ws.Range("I5").Value = 111.333336
$ws.Range("J5").Value = 339
$ws.Range("K5").Value = 111.333336
$ws.Range("L5").Value = 339
$ws.Range("M5").Value = 0.6666639999999973
$ws.Range("N5").Value = -563
$ws.Range("H100").Value = 79700
$ws.Range("J100").Value = 79700
$ws.Range("L100").Value = 79700
$ws.Range("N100").Value = -81864

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6414487.5
$ws.Range("I5").Value = 363.27274
$ws.Range("J5").Value = 41692172
$ws.Range("K5").Value = 1089.81822
$ws.Range("L5").Value = 125076516
$ws.Range("M5").Value = -977.8182200000001
$ws.Range("N5").Value = -125076740
$ws.Range("H12").Value = 41666876
$ws.Range("I12").Value = 142857310
$ws.Range("J12").Value = 219.82353
$ws.Range("K12").Value = 428571930
$ws.Range("L12").Value = 659.47059
$ws.Range("M12").Value = -428571757
$ws.Range("N12").Value = -1005.47059
$ws.Range("H15").Value = 605.3684
$ws.Range("I15").Value = 213.33333
$ws.Range("J15").Value = 786.3077
$ws.Range("K15").Value = 639.99999
$ws.Range("L15").Value = 2358.9231
$ws.Range("M15").Value = -499.99999
$ws.Range("N15").Value = -2638.9231
$ws.Range("H32").Value = 2071.4285
$ws.Range("I32").Value = 50
$ws.Range("J32").Value = 2226.923
$ws.Range("K32").Value = 150
$ws.Range("L32").Value = 6680.768999999999
$ws.Range("M32").Value = 133
$ws.Range("N32").Value = -7246.768999999999
$ws.Range("H68").Value = 3090.9443
$ws.Range("I68").Value = 1510.878
$ws.Range("J68").Value = 5180.7095
$ws.Range("K68").Value = 4532.634
$ws.Range("L68").Value = 15542.1285
$ws.Range("M68").Value = -3721.634
$ws.Range("N68").Value = -17164.1285
$ws.Range("H71").Value = 3090.9443
$ws.Range("I71").Value = 1510.878
$ws.Range("J71").Value = 5180.7095
$ws.Range("K71").Value = 13597.902
$ws.Range("L71").Value = 46626.3855
$ws.Range("M71").Value = -9541.902
$ws.Range("N71").Value = -54738.3855
$ws.Range("H82").Value = 4401.857
$ws.Range("I82").Value = 2271
$ws.Range("K82").Value = 6813
$ws.Range("M82").Value = -6407
$ws.Range("H85").Value = 4401.857
$ws.Range("I85").Value = 2271
$ws.Range("K85").Value = 6813
$ws.Range("M85").Value = -5409
$ws.Range("H107").Value = 579.1548
$ws.Range("J107").Value = 1281.091
$ws.Range("L107").Value = 3843.273
$ws.Range("N107").Value = -7683.272999999999
$ws.Range("H135").Value = 6414487.5
$ws.Range("I135").Value = 363.27274
$ws.Range("J135").Value = 41692172
$ws.Range("K135").Value = 3269.45466
$ws.Range("L135").Value = 375229548
$ws.Range("M135").Value = -734.4546599999999
$ws.Range("N135").Value = -375234618

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H80").Value = 6456.7856
$ws.Range("I80").Value = 9256.286
$ws.Range("K80").Value = 9256.286
$ws.Range("M80").Value = -8258.286
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 6456.7856
$ws.Range("I83").Value = 9256.286
$ws.Range("K83").Value = 46281.43
$ws.Range("M83").Value = -41289.43
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 1949.75
$ws.Range("I17").Value = 933.3333
$ws.Range("J17").Value = 4999
$ws.Range("K17").Value = 933.3333
$ws.Range("L17").Value = 4999
$ws.Range("M17").Value = -763.3333
$ws.Range("N17").Value = -5339

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 150
$ws.Range("I17").Value = 150
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 150
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 22
$ws.Range("N17").ClearContents()
$ws.Range("H39").Value = 10049
$ws.Range("J39").Value = 10049
$ws.Range("L39").Value = 10049
$ws.Range("N39").Value = -10875
$ws.Range("H122").Value = 2334.8333
$ws.Range("I122").Value = 2026
$ws.Range("J122").Value = 2952.5
$ws.Range("K122").Value = 6078
$ws.Range("L122").Value = 8857.5
$ws.Range("M122").Value = -3628
$ws.Range("N122").Value = -13757.5

Write-Output "Applied 187 cell updates across 8 sheets."